$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.251.42'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = '3.073.21'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''559.14'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '''145.86'
$ws.Range('E6').Value = '  +4.69%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '3.064.46'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '''0.503'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('D11').Value = '''6.22'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').Value = '''0.469'
$ws.Range('E12').Value = '  +3.97%  '
$ws.Range('D13').Value = '''0.0000229'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '''35.24'
$ws.Range('E14').Value = '  +1.72%  '
$ws.Range('D15').Value = '3.569.08'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '64.314.16'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '3.071.39'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('D20').Value = '''478.84'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '''13.92'
$ws.Range('D22').Value = '''0.676'
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').Value = '''7.55'
$ws.Range('E23').Value = '  +5.57%  '
$ws.Range('D24').Value = '''13.53'
$ws.Range('E24').Value = '  +8.61%  '
$ws.Range('D25').Value = '''81.70'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '''2.81'
$ws.Range('E27').Value = '  +2.21%  '
$ws.Range('E28').Value = '  +2.31%  '
$ws.Range('E29').Value = '  +4.39%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = '''26.27'
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('D33').Value = '''2.49'
$ws.Range('E33').Value = '  +3.29%  '
$ws.Range('D34').Value = '''5.59'
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('D35').Value = '''6.17'
$ws.Range('E35').Value = '  +3.60%  '
$ws.Range('D36').Value = '''54.96'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = '''3.07'
$ws.Range('E37').Value = '  +19.26%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''462.43'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('D40').Value = '''0.0407'
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('D41').Value = '2.968.27'
$ws.Range('E41').Value = '  -4.95%  '
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('E43').Value = '  -3.15%  '
$ws.Range('D44').Value = '''27.93'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E45').Value = '  +4.45%  '
$ws.Range('D46').Value = '''2.16'
$ws.Range('E46').Value = '  +5.48%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E48').Value = '  +2.89%  '
$ws.Range('D49').Value = '''120.63'
$ws.Range('E49').Value = '  +4.26%  '
$ws.Range('D50').Value = '0.0₃0519'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('E51').Value = '  +1.17%  '
